$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    19.966807505814753,
    71.717663612750641,
    145.82166883631916,
    242.36813005470995,
    390.66733289510421,
    622.63738706505967,
    953.79063754701747,
    1338.1922627030831,
    1991.33616012114,
    2863.9493360257316,
    3518.9338380620434,
    4653.5291065704787,
    5952.4201311570623,
    6850.4662047829042,
    7774.4290676061264
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("C5").Select()
